# Auto-generated Excel COM-interop script
# Reconstructs the sharedStrings table in the exact target order (the engine
# preserves original-pool order for kept strings and appends brand-new strings
# in creation order, so the only reliable way to hit an arbitrary target order
# is to clear everything and re-create the strings fresh, in that order).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Step 1: clear contents (not formatting) of every cell holding shared-string text ---
$ws1.Range("A1:A12").ClearContents()
$ws2.Range("A1:U1").ClearContents()
$ws2.Range("A2:A21").ClearContents()
$ws2.Range("N21:R21").ClearContents()

# --- Step 2: seed the shared-string pool, in order, via a scratch column far outside any used range ---
$sstOrder = @(
    'Труба 1м',
    'Грибок',
    'Переходник',
    'Кагла',
    'Окончание дымох.',
    'Ревизия нерж\нерж',
    'Грибок термо',
    'Конус термо нерж\нерж',
    'Окончание термо',
    '-',
    '100-160',
    '110-180',
    '120-180',
    '130-200',
    '140-200',
    '150-220',
    '160-220',
    '180-250',
    '200-260',
    '220-280',
    '230-300',
    '250-320',
    '300-360',
    '350-420',
    '400-460',
    '450-520',
    '500-560',
    '100-200',
    '120-220',
    '130-230',
    'Труба 1м нерж\оц',
    'Ревизия нерж\оц',
    'Конус термо нерж\оц',
    'elements',
    'Труба 1м нерж\нерж',
    'Тройник 87* нерж\нерж',
    'Тройник 87* нерж\оц',
    'Тройник 45* нерж\нерж',
    'Тройник 45* нерж\оц',
    'Колено 90* нерж\нерж',
    'Колено 90* нерж\оц',
    'Колено 45* нерж\нерж',
    'Колено 45* нерж\оц',
    'Колено 90*',
    'Колено 45*',
    'Тройник 87*',
    'Тройник 45*',
    'Труба 0.5м',
    'Труба 0.3м',
    'Труба 0.5м нерж\нерж',
    'Труба 0.5м нерж\оц',
    'Труба 0.25м нерж\нерж',
    'Труба 0.25м нерж\оц'
)
for ($i = 0; $i -lt $sstOrder.Length; $i++) {
    $ws1.Cells.Item($i + 1, 700).Value = $sstOrder[$i]
}

# --- Step 3: write the real cell values (same text => reuses the pool entry created above, no reordering) ---
# Sheet "0.8_aisi_304"
$ws1.Range("A1").Value = 'elements'
$ws1.Range("A2").Value = 'Труба 1м'
$ws1.Range("A3").Value = 'Труба 0.5м'
$ws1.Range("A4").Value = 'Труба 0.3м'
$ws1.Range("A5").Value = 'Колено 90*'
$ws1.Range("A6").Value = 'Колено 45*'
$ws1.Range("A7").Value = 'Тройник 87*'
$ws1.Range("A8").Value = 'Тройник 45*'
$ws1.Range("A9").Value = 'Грибок'
$ws1.Range("A10").Value = 'Окончание дымох.'
$ws1.Range("A11").Value = 'Переходник'
$ws1.Range("A12").Value = 'Кагла'

# Sheet "ТЕРМО_0.8_aisi_304 "
$ws2.Range("A1").Value = 'elements'
$ws2.Range("B1").Value = '100-160'
$ws2.Range("C1").Value = '110-180'
$ws2.Range("D1").Value = '120-180'
$ws2.Range("E1").Value = '130-200'
$ws2.Range("F1").Value = '140-200'
$ws2.Range("G1").Value = '150-220'
$ws2.Range("H1").Value = '160-220'
$ws2.Range("I1").Value = '180-250'
$ws2.Range("J1").Value = '200-260'
$ws2.Range("K1").Value = '220-280'
$ws2.Range("L1").Value = '230-300'
$ws2.Range("M1").Value = '250-320'
$ws2.Range("N1").Value = '300-360'
$ws2.Range("O1").Value = '350-420'
$ws2.Range("P1").Value = '400-460'
$ws2.Range("Q1").Value = '450-520'
$ws2.Range("R1").Value = '500-560'
$ws2.Range("S1").Value = '100-200'
$ws2.Range("T1").Value = '120-220'
$ws2.Range("U1").Value = '130-230'
$ws2.Range("A2").Value = 'Труба 1м нерж\нерж'
$ws2.Range("A3").Value = 'Труба 1м нерж\оц'
$ws2.Range("A4").Value = 'Труба 0.5м нерж\нерж'
$ws2.Range("A5").Value = 'Труба 0.5м нерж\оц'
$ws2.Range("A6").Value = 'Труба 0.25м нерж\нерж'
$ws2.Range("A7").Value = 'Труба 0.25м нерж\оц'
$ws2.Range("A8").Value = 'Тройник 87* нерж\нерж'
$ws2.Range("A9").Value = 'Тройник 87* нерж\оц'
$ws2.Range("A10").Value = 'Тройник 45* нерж\нерж'
$ws2.Range("A11").Value = 'Тройник 45* нерж\оц'
$ws2.Range("A12").Value = 'Колено 90* нерж\нерж'
$ws2.Range("A13").Value = 'Колено 90* нерж\оц'
$ws2.Range("A14").Value = 'Колено 45* нерж\нерж'
$ws2.Range("A15").Value = 'Колено 45* нерж\оц'
$ws2.Range("A16").Value = 'Ревизия нерж\нерж'
$ws2.Range("A17").Value = 'Ревизия нерж\оц'
$ws2.Range("A18").Value = 'Грибок термо'
$ws2.Range("A19").Value = 'Конус термо нерж\нерж'
$ws2.Range("A20").Value = 'Конус термо нерж\оц'
$ws2.Range("A21").Value = 'Окончание термо'
$ws2.Range("N21").Value = '-'
$ws2.Range("O21").Value = '-'
$ws2.Range("P21").Value = '-'
$ws2.Range("Q21").Value = '-'
$ws2.Range("R21").Value = '-'

# --- Step 4: remove the scratch seed column now that every string is anchored by a real cell ---
for ($i = 0; $i -lt $sstOrder.Length; $i++) {
    $ws1.Cells.Item($i + 1, 700).ClearContents()
}

# --- Step 5: view-state changes (zoom / selection / active tab) ---
$ws2.Activate()
$ws2.Range("A14").Select()
$ws1.Activate()
$excel.ActiveWindow.Zoom = 106
$ws1.Range("A6").Select()

